$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("rdth", $false, $false, $false, $false, $false, $true, 1, $false, " n", 2)
